# Apply cryptos.xlsx data refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.609.49'
$ws.Range("E2").Value = '  -1.26%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.782.28'
$ws.Range("E3").Value = '  +0.82%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.27'
$ws.Range("E5").Value = '  -0.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.70'
$ws.Range("E6").Value = '  -0.30%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.781.36'
$ws.Range("E7").Value = '  +0.83%  '

# Row 9
$ws.Range("E9").Value = '  -0.41%  '

# Row 10
$ws.Range("E10").Value = '  +0.52%  '

# Row 11
$ws.Range("E11").Value = '  -2.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  +0.23%  '

# Row 13
$ws.Range("E13").Value = '  -0.86%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.38'
$ws.Range("E14").Value = '  +0.64%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.418.56'
$ws.Range("E15").Value = '  +0.87%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.806.41'
$ws.Range("E16").Value = '  +1.83%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.43'
$ws.Range("E17").Value = '  +3.01%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.564.93'
$ws.Range("E18").Value = '  -1.34%  '

# Row 19
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("E20").Value = '  -0.25%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").Value = '  -4.95%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.31'
$ws.Range("E22").Value = '  -2.48%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.696'
$ws.Range("E23").Value = '  -0.45%  '

# Row 24
$ws.Range("E24").Value = '  +7.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.55'
$ws.Range("E25").Value = '  -0.77%  '

# Row 26
$ws.Range("E26").Value = '  -1.13%  '

# Row 27
$ws.Range("E27").Value = '  -2.54%  '

# Row 28
$ws.Range("E28").Value = '  -0.71%  '

# Row 29
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("E30").Value = '  +0.32%  '

# Row 31
$ws.Range("E31").Value = '  -0.05%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.81'
$ws.Range("E32").Value = '  -0.29%  '

# Row 33
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("E33").Value = '  +0.58%  '

# Row 34
$ws.Range("E34").Value = '  -0.20%  '

# Row 35
$ws.Range("E35").Value = '  -0.28%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.737.11'
$ws.Range("E36").Value = '  +0.83%  '

# Row 37
$ws.Range("E37").Value = '  -1.21%  '

# Row 38
$ws.Range("E38").Value = '  -2.02%  '

# Row 39
$ws.Range("E39").Value = '  -0.49%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("E40").Value = '  -0.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.76'
$ws.Range("E41").Value = '  -0.73%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.37'
$ws.Range("E44").Value = '  +5.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.299'
$ws.Range("E45").Value = '  -1.98%  '

# Row 46
$ws.Range("E46").Value = '  +2.79%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.34'
$ws.Range("E47").Value = '  -2.96%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '148.81'
$ws.Range("E48").Value = '  +1.70%  '

# Row 49
$ws.Range("E49").Value = '  -4.74%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '389.49'
$ws.Range("E50").Value = '  -0.52%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.65'
$ws.Range("E51").Value = '  +2.00%  '
